# Adapt column header formatting to respective input file names (#7)
#   *_old  -> *_FV2310
#   *_new  -> *_FV2404
# Then turn the sheet's used range into an Excel Table (ListObject) and
# freeze the header row, matching the target workbook layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header cells (A1:U1) ---------------------------------------
$headerMap = @{
    "Segmentname_old"         = "Segmentname_FV2310"
    "Segmentgruppe_old"       = "Segmentgruppe_FV2310"
    "Segment_old"             = "Segment_FV2310"
    "Datenelement_old"        = "Datenelement_FV2310"
    "Segment ID_old"          = "Segment ID_FV2310"
    "Code_old"                = "Code_FV2310"
    "Qualifier_old"           = "Qualifier_FV2310"
    "Beschreibung_old"        = "Beschreibung_FV2310"
    "Bedingungsausdruck_old"  = "Bedingungsausdruck_FV2310"
    "Bedingung_old"           = "Bedingung_FV2310"
    "Segmentname_new"         = "Segmentname_FV2404"
    "Segmentgruppe_new"       = "Segmentgruppe_FV2404"
    "Segment_new"             = "Segment_FV2404"
    "Datenelement_new"        = "Datenelement_FV2404"
    "Segment ID_new"          = "Segment ID_FV2404"
    "Code_new"                = "Code_FV2404"
    "Qualifier_new"           = "Qualifier_FV2404"
    "Beschreibung_new"        = "Beschreibung_FV2404"
    "Bedingungsausdruck_new"  = "Bedingungsausdruck_FV2404"
    "Bedingung_new"           = "Bedingung_FV2404"
}

$headerRange = $ws.Range("A1:U1")
$colCount = $headerRange.Columns.Count
for ($i = 1; $i -le $colCount; $i++) {
    $cell = $headerRange.Cells.Item(1, $i)
    $cur = $cell.Value()
    if ($headerMap.ContainsKey($cur)) {
        $cell.Value = $headerMap[$cur]
    }
}

# --- 2) Turn the used range into an Excel Table (ListObject) -------------
# Temporarily strip the header row's manual formatting so the engine does
# not capture it as a bespoke header-row dxf on the new table.
$dataRange = $ws.Range("A1:U85")
$headerRange.Style = "Normal"

$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Drop the implicit table style so no extra style/dxf metadata is added.
$tbl.TableStyle = ""

# Re-apply the header row's original look (bold, centered, wrapped, filled,
# thin-bordered) now that the table exists.
$headerRange.Borders.LineStyle = 1
$headerRange.WrapText = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.Interior.Color = 14277081
$headerRange.Font.Bold = $true

# --- 3) Freeze the header row ---------------------------------------------
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
